$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell H1 ("Save"), matching the formatting used by the other
# header cells (B1:G1) by copying G1's format (bold, centered, bordered) onto it.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add values for the new "Save" column
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
